$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the "regulator" data row from the "appointment" table
#    (the first table in the document - it is the last row there).
#    Note: Range.Text on a table cell carries trailing control chars
#    (cell-end mark \r\a), so compare with StartsWith rather than -eq.
# ------------------------------------------------------------------
$appointmentTable = $d.Tables.Item(1)
$lastRow = $appointmentTable.Rows.Item($appointmentTable.Rows.Count)
if ($lastRow.Cells.Item(1).Range.Text.StartsWith("regulator")) {
    $lastRow.Delete()
}

# ------------------------------------------------------------------
# 2) Remove the whole "regulator" type section at the end of the
#    document: its Heading 1 paragraph plus the definition table
#    that follows it.
#    (Re-deriving the Paragraphs collection from a fresh Range after
#    the row delete above keeps indices accurate.)
# ------------------------------------------------------------------
$bodyRange = $d.Content
$paraCount = $bodyRange.Paragraphs.Count
$headingIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $bodyRange.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("regulator") -and $para.Range.Text.Length -eq 10 -and $para.Style.NameLocal -eq "Heading 1") {
        $headingIndex = $i
    }
}

if ($headingIndex -ne -1) {
    $headingPara = $bodyRange.Paragraphs.Item($headingIndex)
    $headingPara.Range.Delete()
}

$regulatorTable = $d.Tables.Item($d.Tables.Count)
if ($regulatorTable.Rows.Item(1).Cells.Item(1).Range.Text.StartsWith("Nom de balise") -and `
    $regulatorTable.Rows.Item(2).Cells.Item(1).Range.Text.StartsWith("regulatorId")) {
    $regulatorTable.Delete()
}
